$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New manufacturer rows
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "AMD"

$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "Intel"

$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "Nvidia"

# Update selection to match target state
$ws.Range("E9").Select()
